$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "m3_s"
$ws.Range("C1").Value = "m3_m"
$ws.Range("D1").Value = "m3_l"
$ws.Range("E1").Value = "m_rpt_s"
$ws.Range("F1").Value = "m_rpt_m"
$ws.Range("G1").Value = "m_rpt_l"

$ws.Range("M7").Select()
